$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. B51: "1.3.0" -> "1.4.0"
$ws.Range("B51").Value = "1.4.0"

# 2. I52: reorder the comma separated ValueSet list (same set of items)
$ws.Range("I52").Value = "ValueSet/us-core-diagnosticreport-category,ValueSet/c80-doc-typecodes,ValueSet/us-core-documentreference-type,ValueSet/2.16.840.1.113883.3.88.12.80.62,ValueSet/document-classcodes,ValueSet/2.16.840.1.113762.1.4.1267.13,ValueSet/us-core-clinical-note-type,ValueSet/us-core-goal-description,ValueSet/us-core-laboratory-test-codes,ValueSet/us-core-diagnosticreport-report-and-note-codes,ValueSet/care-team-category,ValueSet/2.16.840.1.113883.11.20.9.69.4,ValueSet/observation-codes,ValueSet/us-core-procedure-code"

# 3. H53: reorder the comma separated ValueSet list
$ws.Range("H53").Value = "ValueSet/servicerequest-category,ValueSet/specimen-container-type,ValueSet/2.16.840.1.113762.1.4.1099.54,ValueSet/procedure-followup,ValueSet/us-core-goal-description,ValueSet/us-core-condition-code-current,ValueSet/specimen-collection-method,ValueSet/medication-form-codes,ValueSet/body-site,ValueSet/manifestation-or-symptom,ValueSet/procedure-code,ValueSet/goal-start-event,ValueSet/medication-as-needed-reason,ValueSet/2.16.840.1.113762.1.4.1166.22,ValueSet/c80-facilitycodes,ValueSet/immunization-reason,ValueSet/2.16.840.1.113762.1.4.1267.3,ValueSet/substance-code,ValueSet/us-core-procedure-code,ValueSet/observation-methods,ValueSet/care-plan-category,ValueSet/questionnaire-answers,ValueSet/servicerequest-orderdetail,ValueSet/clinical-findings,ValueSet/condition-stage,ValueSet/immunization-target-disease,ValueSet/encounter-reason,ValueSet/additional-instruction-codes,ValueSet/participant-role,ValueSet/approach-site-codes,ValueSet/device-action,ValueSet/us-core-condition-code,ValueSet/medication-codes,ValueSet/device-kind,ValueSet/performer-role,ValueSet/c80-practice-codes,ValueSet/referencerange-appliesto,ValueSet/us-core-servicerequest-category,ValueSet/route-codes,ValueSet/care-plan-activity-outcome,ValueSet/procedure-outcome,ValueSet/procedure-category,ValueSet/2.16.840.1.113762.1.4.1021.121,ValueSet/condition-code,ValueSet/condition-stage-type,ValueSet/procedure-not-performed-reason,ValueSet/immunization-status-reason,ValueSet/administration-method-codes,ValueSet/condition-severity"

# 4. H125: reorder the comma separated ValueSet list
$ws.Range("H125").Value = "ValueSet/ucum-bodyweight,ValueSet/ucum-bodytemp,ValueSet/ucum-common,ValueSet/ucum-vitals-common,ValueSet/ucum-bodylength,ValueSet/v3-UnitsOfMeasureCaseSensitive"

# 5. Row 130: replace HSLOC entry with CMS Place of Service Codes entry
$ws.Range("A130").Value = "https://www.cms.gov/Medicare/Coding/place-of-service-codes/Place_of_Service_Code_Set"
$ws.Range("B130").Value = 20240502
$ws.Range("C130").Value = "active"
$ws.Range("D130").Value = "2.16.840.1.113883.15.5"
$ws.Range("E130").Value = "CMSPlaceofServiceCodes"
$ws.Range("F130").Value = "CMS Place of Service Codes (POS)"
$ws.Range("G130").Value = "Place of Service Codes are two-digit codes placed on health care professional claims to indicate the setting in which a service was provided. The Centers for Medicare & Medicaid Services (CMS) maintain POS codes used throughout the health care industry.    This code set is required for use in the implementation guide adopted as the national standard for electronic transmission of professional health care claims under the provisions of the Health Insurance Portability and Accountability Act of 1996 (HIPAA). HIPAA directed the Secretary of HHS to adopt national standards for electronic transactions. These standard transactions require all health plans and providers to use standard code sets to populate data elements in each transaction. The Transaction and Code Set Rule adopted the ASC X12N-837 Health Care Claim: Professional, volumes 1 and 2, as the standard for electronic submission of professional claims. This standard names the POS code set currently maintained by CMS as the code set to be used for describing sites of service in such claims. POS information is often needed to determine the acceptability of direct billing of Medicare, Medicaid and private insurance services provided by a given provider."
$ws.Range("H130").Value = "ValueSet/CMSPlaceOfServiceCodes"

# 6. Insert a new row at 131, pushing the old "urn:ietf:bcp:47" row down to 132,
#    then populate the new row 131 with the USPS State Codes entry.
$ws.Rows.Item(131).Insert()

$ws.Range("A131").Value = "https://www.usps.com/"
$ws.Range("B131").Value = "1.0.0"
$ws.Range("C131").Value = "active"
$ws.Range("D131").Value = ""
$ws.Range("E131").Value = "USPSStateCodes"
$ws.Range("F131").Value = "USPS State Codes"
$ws.Range("G131").Value = "State Codes as defined by US Postal Service    For more information, see [https://pe.usps.com/text/pub28/28apb.htm](https://pe.usps.com/text/pub28/28apb.htm)"
$ws.Range("H131").Value = "ValueSet/USPS-State"
$ws.Range("I131").Value = ""
$ws.Range("J131").Value = ""
